# Weekly update: insert the newest Fruta/Mango price record at row 382 of
# the "Feria Lagunitas de Puerto Montt" price series. Excel's row-insert
# semantics push the existing row 382 (and everything below it) down by one
# row, so the historical rows simply shift from row N to row N+1, and what
# used to be the last row (441) becomes the new last row (442).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 382, shifting rows 382..441
# down to 383..442 (this also grows the sheet's used range to A1:T442).
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with this week's record.
$ws.Range("A382").Value = 4
$ws.Range("B382").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C382").Value = "Los Lagos"
$ws.Range("D382").Value = 45180
$ws.Range("E382").Value = 10
$ws.Range("F382").Value = "Fruta"
$ws.Range("G382").Value = 100108
$ws.Range("H382").Value = "Tropicales y subtropicales"
$ws.Range("I382").Value = 100108002
$ws.Range("J382").Value = "Mango"
$ws.Range("K382").Value = "Sin especificar"
$ws.Range("L382").Value = "Primera"
$ws.Range("M382").Value = 60
$ws.Range("N382").Value = 12000
$ws.Range("O382").Value = 12000
$ws.Range("P382").Value = 12000
$ws.Range("Q382").Value = "$/bandeja 4 kilos"
$ws.Range("R382").Value = "Brasil"
$ws.Range("S382").Value = 3000
$ws.Range("T382").Value = 4
